$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title "F050A" -> "F050I": the big heading on the first paragraph
#    is split into three separate runs ("F", "050", "A"). We only want
#    to retarget the lone "A" run to "I", leaving "F"/"050" untouched,
#    so the Find is scoped to paragraph 1's Range (ReplaceOne = 1).
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.Find.Execute("A", $true, $false, $false, $false, $false, `
    $true, 1, $false, "I", 1) | Out-Null

# ------------------------------------------------------------------
# 2) "F050A" -> "F050I" in the "Board ID" table cell, and
#    "German" -> "Italian" in the "Language" table cell.
#    Both cells hold their text in a single run already, but we still
#    scope each Find to its own paragraph so only that one cell is hit.
# ------------------------------------------------------------------
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text

    if ($text -like "F050A*") {
        $r = $para.Range
        $r.Find.Execute("F050A", $true, $false, $false, $false, $false, `
            $true, 1, $false, "F050I", 1) | Out-Null
    }
    elseif ($text -like "German*") {
        $r = $para.Range
        $r.Find.Execute("German", $true, $false, $false, $false, $false, `
            $true, 1, $false, "Italian", 1) | Out-Null
    }
}

# ------------------------------------------------------------------
# 3) Fill in the empty run right below "Notes: " with three lines of
#    text separated by manual line breaks (<w:br/>), keeping all three
#    <w:t> pieces inside the single existing run/rPr (a literal
#    vertical-tab char, Chr(11), is how Word COM encodes a manual line
#    break without splitting the run).
# ------------------------------------------------------------------
$notesLabelIndex = -1
for ($i = 1; $i -le $n; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "Notes:*") {
        $notesLabelIndex = $i
        break
    }
}

if ($notesLabelIndex -gt 0) {
    $notesPara = $d.Paragraphs($notesLabelIndex + 1)
    $notesRange = $notesPara.Range
    $lineBreak = [char]11
    $notesText = "scratches on bottom" + $lineBreak + `
        "exposed copper on bottom" + $lineBreak + `
        "large scratches on bottom"
    $notesRange.Text = $notesText
}
